# Append new quiz-result rows (Paul & Katy's data) to the "Worksheet" sheet.
# Columns: A=Student Id, B=Trial, C=Correct, D=Elapsed Time, E=Date
# Note: values in columns C ("false") and E (date-looking text) must stay as
# literal text, so a leading apostrophe is used to stop Excel auto-typing them
# as Boolean/Date, and ClearFormats() strips the resulting quote-prefix style
# so the cell keeps the workbook's default (unstyled) formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("pstdenis@stonybrook.edu",   1, 248, "2019-12-24"),
    @("pstdenis@stonybrook.edu",   2, 298, "2019-12-24"),
    @("pstdenis@stonybrook.edu",   1, 13,  "2019-12-24"),
    @("pstdenis@stonybrook.edu",   2, 15,  "2019-12-24"),
    @("asklyarova@stonybrook.edu", 1, 113, "2019-12-26"),
    @("asklyarova@stonybrook.edu", 2, 123, "2019-12-26"),
    @("ikleiman@stonybrook.edu",   2, 9,   "2019-12-30")
)

$startRow = 14
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $row = $startRow + $i
    $vals = $newRows[$i]

    $ws.Cells.Item($row, 1).Value = $vals[0]

    $ws.Cells.Item($row, 2).Value = $vals[1]

    $cCell = $ws.Cells.Item($row, 3)
    $cCell.Value = "'false"
    $cCell.ClearFormats()

    $ws.Cells.Item($row, 4).Value = $vals[2]

    $eCell = $ws.Cells.Item($row, 5)
    $eCell.Value = "'" + $vals[3]
    $eCell.ClearFormats()
}
